# Apply cryptos list update (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "68.447.32"
Set-TextCell $ws.Range("E2") "  -0.66%  "

Set-TextCell $ws.Range("D3") "3.853.99"

Set-TextCell $ws.Range("E4") "  +0.08%  "

Set-TextCell $ws.Range("D5") "520.78"
Set-TextCell $ws.Range("E5") "  +6.68%  "

Set-TextCell $ws.Range("D6") "140.44"
Set-TextCell $ws.Range("E6") "  -4.09%  "

Set-TextCell $ws.Range("E7") "  -2.45%  "

Set-TextCell $ws.Range("E8") "  +0.13%  "

Set-TextCell $ws.Range("D9") "0.711"
Set-TextCell $ws.Range("E9") "  -3.30%  "

Set-TextCell $ws.Range("E10") "  -5.65%  "

Set-TextCell $ws.Range("E11") "  -7.96%  "

Set-TextCell $ws.Range("D12") "41.35"
Set-TextCell $ws.Range("E12") "  -3.76%  "

Set-TextCell $ws.Range("D13") "10.30"
Set-TextCell $ws.Range("E13") "  -1.50%  "

Set-TextCell $ws.Range("D14") "4.470.26"
Set-TextCell $ws.Range("E14") "  -1.92%  "

Set-TextCell $ws.Range("D15") "21.26"
Set-TextCell $ws.Range("E15") "  +6.69%  "

Set-TextCell $ws.Range("D16") "3.840.88"
Set-TextCell $ws.Range("E16") "  -2.37%  "

Set-TextCell $ws.Range("D17") "14.07"
Set-TextCell $ws.Range("E17") "  -1.23%  "

Set-TextCell $ws.Range("E18") "  -2.10%  "

Set-TextCell $ws.Range("D19") "1.19"
Set-TextCell $ws.Range("E19") "  +2.41%  "

Set-TextCell $ws.Range("D20") "68.466.10"
Set-TextCell $ws.Range("E20") "  -0.76%  "

Set-TextCell $ws.Range("D21") "415.44"
Set-TextCell $ws.Range("E21") "  -4.89%  "

Set-TextCell $ws.Range("D22") "3.47"
Set-TextCell $ws.Range("E22") "  +0.35%  "

Set-TextCell $ws.Range("E23") "  -3.80%  "

Set-TextCell $ws.Range("D24") "86.57"
Set-TextCell $ws.Range("E24") "  -3.17%  "

Set-TextCell $ws.Range("E25") "  +6.26%  "

Set-TextCell $ws.Range("D26") "11.51"
Set-TextCell $ws.Range("E26") "  -7.40%  "

Set-TextCell $ws.Range("D27") "10.51"
Set-TextCell $ws.Range("E27") "  -5.54%  "

Set-TextCell $ws.Range("D28") "35.33"
Set-TextCell $ws.Range("E28") "  -4.86%  "

Set-TextCell $ws.Range("D29") "13.17"
Set-TextCell $ws.Range("E29") "  -2.60%  "

Set-TextCell $ws.Range("E30") "  -4.62%  "

Set-TextCell $ws.Range("E31") "  -6.31%  "

Set-TextCell $ws.Range("B32") "Toncoin"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws.Range("D32") "2.78"
Set-TextCell $ws.Range("E32") "  -3.85%  "

Set-TextCell $ws.Range("B33") "NEARProtocol"
Set-TextCell $ws.Range("C33") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D33") "6.64"
Set-TextCell $ws.Range("E33") "  +9.68%  "

Set-TextCell $ws.Range("D34") "65.88"
Set-TextCell $ws.Range("E34") "  +6.56%  "

Set-TextCell $ws.Range("D35") "0.455"
Set-TextCell $ws.Range("E35") "  -5.49%  "

Set-TextCell $ws.Range("D36") "39.53"
Set-TextCell $ws.Range("E36") "  -3.03%  "

Set-TextCell $ws.Range("E37") "  +14.80%  "

Set-TextCell $ws.Range("D38") "0.0₃0825"
Set-TextCell $ws.Range("E38") "  -7.58%  "

Set-TextCell $ws.Range("D39") "0.149"
Set-TextCell $ws.Range("E39") "  -0.86%  "

Set-TextCell $ws.Range("E40") "  -0.06%  "

Set-TextCell $ws.Range("E41") "  -0.17%  "

Set-TextCell $ws.Range("D42") "0.0472"
Set-TextCell $ws.Range("E42") "  -3.63%  "

Set-TextCell $ws.Range("D43") "3.10"
Set-TextCell $ws.Range("E43") "  +3.22%  "

Set-TextCell $ws.Range("D44") "2.79"
Set-TextCell $ws.Range("E44") "  -5.71%  "

Set-TextCell $ws.Range("D45") "3.38"
Set-TextCell $ws.Range("E45") "  +0.69%  "

Set-TextCell $ws.Range("D46") "0.140"
Set-TextCell $ws.Range("E46") "  -2.94%  "

Set-TextCell $ws.Range("D47") "3.00"
Set-TextCell $ws.Range("E47") "  -0.06%  "

Set-TextCell $ws.Range("D48") "0.000266"
Set-TextCell $ws.Range("E48") "  +11.81%  "

Set-TextCell $ws.Range("D49") "143.44"
Set-TextCell $ws.Range("E49") "  +0.29%  "

Set-TextCell $ws.Range("D50") "3.25"
Set-TextCell $ws.Range("E50") "  -4.15%  "

Set-TextCell $ws.Range("D51") "0.0₆0337"
Set-TextCell $ws.Range("E51") "  -6.52%  "

